$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.315.25"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "1.869.29"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4705"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.52%  "

$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("E9").Value = "  +0.50%  "

$ws.Range("E10").Value = "  -3.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08012"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.03"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.86%  "

$ws.Range("D13").Value = "1.869.83"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.119"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6848"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "269.97"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.49%  "

$ws.Range("D17").Value = "30.309.11"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.02"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007632"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.98%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").Value = "2.114.00"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.292"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.20%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.221"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.448"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.75"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.92"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.950"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.369"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.93%  "

$ws.Range("E30").Value = "  +0.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.370"
$ws.Range("D31").ClearFormats()

$ws.Range("E32").Value = "  -1.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.073"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04712"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.85%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6998"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.699"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01867"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.633"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.296"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.06"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -6.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.961"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8426"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4171"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9995"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.89"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.050"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.57%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.144"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "912.19"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.50"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05699"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.94%  "
